# updated main GSC export data
#
# The "Chart" sheet holds a rolling daily export: column A = date (stored as
# plain text, e.g. "2025-11-02"), column B = Invalid count, column C = Valid
# count - one row per day. The export window rolls forward by one day: the
# oldest day (first data row) is dropped, every remaining day shifts up one
# row, and a new trailing day is appended at the end (carrying forward the
# last known Valid count, since the newest day's GSC data is typically not
# fully processed/crawled yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$firstDataRow = 2
$lastDataRow = $ws.UsedRange.Rows.Count

# A cell whose number format ("General") is never touched by this edit - used
# below to strip the "Text" number format Excel applies when a date-looking
# string is assigned, so the date cells stay plain text/General like the
# original export instead of turning into real date serials.
$formatDonor = $ws.Cells.Item($firstDataRow, 2)

# Snapshot the current (pre-shift) column A/C values before writing anything,
# since the shift reads row+1's value while writing row r.
$dates = @{}
$values = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $dates[$r] = $ws.Cells.Item($r, 1).Value2
    $values[$r] = $ws.Cells.Item($r, 3).Value2
}

# Compute the new trailing day from the last known date (+1 day); the new
# day's Valid count simply carries forward the prior last value.
$lastDate = [datetime]::ParseExact($dates[$lastDataRow], "yyyy-MM-dd", $null)
$newDateStr = $lastDate.AddDays(1).ToString("yyyy-MM-dd")
$newValue = $values[$lastDataRow]
$dates[$lastDataRow + 1] = $newDateStr
$values[$lastDataRow + 1] = $newValue

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $dateCell = $ws.Cells.Item($r, 1)

    # Force text entry (otherwise Excel auto-converts "yyyy-MM-dd"-looking
    # text into a date value), then immediately restore a plain/General
    # number format so the cell matches the original export's formatting.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dates[$r + 1]
    $formatDonor.Copy()
    $dateCell.PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 3).Value = $values[$r + 1]
}

$excel.CutCopyMode = 0
